$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need NumberFormat = "@"
# set first so Excel keeps them as literal text instead of coercing to a number,
# matching the source data (these are price strings, not numeric values).
$textCells = @("D5", "D8", "D10", "D11", "D20", "D25", "D26", "D28", "D33", "D36", "D40", "D42", "D44", "D51")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = '35.216.79'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = '1.861.11'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E4").Value = '  +0.77%  '
$ws.Range("D5").Value = '239.22'
$ws.Range("E5").Value = '  +3.50%  '
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("E7").Value = '  +0.72%  '
$ws.Range("D8").Value = '42.39'
$ws.Range("E8").Value = '  +6.45%  '
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("D10").Value = '0.0693'
$ws.Range("E10").Value = '  +1.27%  '
$ws.Range("D11").Value = '0.0991'
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").Value = '2.129.47'
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("D14").Value = '1.861.97'
$ws.Range("E14").Value = '  +1.17%  '
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("D17").Value = '35.175.46'
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("E19").Value = '  +1.07%  '
$ws.Range("D20").Value = '241.49'
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("E22").Value = '  +1.18%  '
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("D25").Value = '169.36'
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("D26").Value = '1.88'
$ws.Range("E26").Value = '  +24.19%  '
$ws.Range("E27").Value = '  +3.29%  '
$ws.Range("D28").Value = '17.72'
$ws.Range("E28").Value = '  +1.72%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  +1.71%  '
$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("E32").Value = '  +1.73%  '
$ws.Range("D33").Value = '1.83'
$ws.Range("E33").Value = '  +28.71%  '
$ws.Range("E34").Value = '  +1.88%  '
$ws.Range("E35").Value = '  +9.48%  '
$ws.Range("D36").Value = '0.820'
$ws.Range("E36").Value = '  +17.86%  '
$ws.Range("E37").Value = '  +7.56%  '
$ws.Range("E38").Value = '  +2.95%  '
$ws.Range("D40").Value = '90.18'
$ws.Range("E40").Value = '  -1.20%  '
$ws.Range("D41").Value = '1.348.49'
$ws.Range("E41").Value = '  +0.56%  '
$ws.Range("D42").Value = '0.0599'
$ws.Range("E42").Value = '  +14.73%  '
$ws.Range("E43").Value = '  +3.09%  '
$ws.Range("D44").Value = '2.33'
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("E46").Value = '  +47.72%  '
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("E48").Value = '  +4.66%  '
$ws.Range("D49").Value = '2.043.91'
$ws.Range("E49").Value = '  +1.24%  '
$ws.Range("E50").Value = '  +0.93%  '
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").Value = '1.01'
$ws.Range("E51").Value = '  +0.73%  '
